# KW-RANK B08MPM2BB2-1-100.xlsx 2020-12-05 07:10:44
# Renames the sheet to the new ASIN and replaces the keyword list in column A
# with the new ranked keyword set, then resets the cell selection to A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to match the new ASIN (B08L7BD4C6 -> B08MPM2BB2).
$ws.Name = "B08MPM2BB2"

# Replace the 100 keyword rows in column A with the new keyword list, in order.
$ws.Cells.Item(1, 1).Value = "black top"
$ws.Cells.Item(2, 1).Value = "black bra"
$ws.Cells.Item(3, 1).Value = "wireless bra"
$ws.Cells.Item(4, 1).Value = "sport bra"
$ws.Cells.Item(5, 1).Value = "pink top"
$ws.Cells.Item(6, 1).Value = "yoga bra"
$ws.Cells.Item(7, 1).Value = "yoga top"
$ws.Cells.Item(8, 1).Value = "low back bra"
$ws.Cells.Item(9, 1).Value = "black s"
$ws.Cells.Item(10, 1).Value = "longline bra"
$ws.Cells.Item(11, 1).Value = "padded bra"
$ws.Cells.Item(12, 1).Value = "green bra"
$ws.Cells.Item(13, 1).Value = "strappy bra"
$ws.Cells.Item(14, 1).Value = "bra top"
$ws.Cells.Item(15, 1).Value = "pink bra"
$ws.Cells.Item(16, 1).Value = "pink bralette"
$ws.Cells.Item(17, 1).Value = "workout bra"
$ws.Cells.Item(18, 1).Value = "strappy bralette"
$ws.Cells.Item(19, 1).Value = "black bra top"
$ws.Cells.Item(20, 1).Value = "longline bralette"
$ws.Cells.Item(21, 1).Value = "black cross"
$ws.Cells.Item(22, 1).Value = "bralette top"
$ws.Cells.Item(23, 1).Value = "cute top"
$ws.Cells.Item(24, 1).Value = "criss cross bra"
$ws.Cells.Item(25, 1).Value = "black l"
$ws.Cells.Item(26, 1).Value = "black m"
$ws.Cells.Item(27, 1).Value = "wirefree bra"
$ws.Cells.Item(28, 1).Value = "criss cross top"
$ws.Cells.Item(29, 1).Value = "long line bra"
$ws.Cells.Item(30, 1).Value = "cute bra"
$ws.Cells.Item(31, 1).Value = "criss cross"
$ws.Cells.Item(32, 1).Value = "pink bra top"
$ws.Cells.Item(33, 1).Value = "bralette bra"
$ws.Cells.Item(34, 1).Value = "impact sport"
$ws.Cells.Item(35, 1).Value = "long s"
$ws.Cells.Item(36, 1).Value = "top s"
$ws.Cells.Item(37, 1).Value = "black cup"
$ws.Cells.Item(38, 1).Value = "wireless bralette"
$ws.Cells.Item(39, 1).Value = "women’s longline sports bra wirefree padded medium support yoga bras gym running workout tank tops"
$ws.Cells.Item(40, 1).Value = "cross bra"
$ws.Cells.Item(41, 1).Value = "green l"
$ws.Cells.Item(42, 1).Value = "cup with"
$ws.Cells.Item(43, 1).Value = "sport gym"
$ws.Cells.Item(44, 1).Value = "long bra"
$ws.Cells.Item(45, 1).Value = "yoga gym"
$ws.Cells.Item(46, 1).Value = "low back bralette"
$ws.Cells.Item(47, 1).Value = "longline yoga bra"
$ws.Cells.Item(48, 1).Value = "criss cross bralette"
$ws.Cells.Item(49, 1).Value = "sport elastic"
$ws.Cells.Item(50, 1).Value = "cute back"
$ws.Cells.Item(51, 1).Value = "pink cross"
$ws.Cells.Item(52, 1).Value = "green s"
$ws.Cells.Item(53, 1).Value = "black apparel"
$ws.Cells.Item(54, 1).Value = "black yoga top"
$ws.Cells.Item(55, 1).Value = "green m"
$ws.Cells.Item(56, 1).Value = "yoga bra strappy"
$ws.Cells.Item(57, 1).Value = "cross back bra"
$ws.Cells.Item(58, 1).Value = "cute pink"
$ws.Cells.Item(59, 1).Value = "yoga workout"
$ws.Cells.Item(60, 1).Value = "sport workout"
$ws.Cells.Item(61, 1).Value = "cute bralette"
$ws.Cells.Item(62, 1).Value = "low back top"
$ws.Cells.Item(63, 1).Value = "criss cross back top"
$ws.Cells.Item(64, 1).Value = "strappy top"
$ws.Cells.Item(65, 1).Value = "criss cross front"
$ws.Cells.Item(66, 1).Value = "strappy sport bra"
$ws.Cells.Item(67, 1).Value = "cross front bra"
$ws.Cells.Item(68, 1).Value = "bra elastic"
$ws.Cells.Item(69, 1).Value = "long line"
$ws.Cells.Item(70, 1).Value = "sport top"
$ws.Cells.Item(71, 1).Value = "strappy back bralette"
$ws.Cells.Item(72, 1).Value = "green sport bra"
$ws.Cells.Item(73, 1).Value = "bra cup"
$ws.Cells.Item(74, 1).Value = "wireless sport"
$ws.Cells.Item(75, 1).Value = "top m"
$ws.Cells.Item(76, 1).Value = "strappy back top"
$ws.Cells.Item(77, 1).Value = "gym apparel"
$ws.Cells.Item(78, 1).Value = "gym back"
$ws.Cells.Item(79, 1).Value = "back top"
$ws.Cells.Item(80, 1).Value = "strappy yoga bra"
$ws.Cells.Item(81, 1).Value = "sport cup"
$ws.Cells.Item(82, 1).Value = "strappy back bra"
$ws.Cells.Item(83, 1).Value = "strappy yoga top"
$ws.Cells.Item(84, 1).Value = "criss cross back"
$ws.Cells.Item(85, 1).Value = "long l"
$ws.Cells.Item(86, 1).Value = "low front bra"
$ws.Cells.Item(87, 1).Value = "gym bra"
$ws.Cells.Item(88, 1).Value = "workout sport bra"
$ws.Cells.Item(89, 1).Value = "cross back bralette"
$ws.Cells.Item(90, 1).Value = "pink criss cross top"
$ws.Cells.Item(91, 1).Value = "black strappy top"
$ws.Cells.Item(92, 1).Value = "cross line"
$ws.Cells.Item(93, 1).Value = "strappy workout bra"
$ws.Cells.Item(94, 1).Value = "back workout"
$ws.Cells.Item(95, 1).Value = "yoga bra top"
$ws.Cells.Item(96, 1).Value = "black strappy"
$ws.Cells.Item(97, 1).Value = "black criss cross top"
$ws.Cells.Item(98, 1).Value = "padded yoga bra"
$ws.Cells.Item(99, 1).Value = "pink bralette top"
$ws.Cells.Item(100, 1).Value = "green cross"

# Reset the selection back to A1 (was E15 previously).
[void]$ws.Range("A1").Select()

